# Adds a new "LIS" tracking row (row 22) below the existing "Coin Change"
# row, fixes up the date-format style on A21, and scrolls the sheet view so
# the newly added row is visible (topLeftCell = A10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Same date serial (46081 -> 2026-02-28) as the previous row - solved same day.
$lisDate = $ws.Range("A21").Value

# Normalize A21's number format style to match the rest of the date column.
$ws.Range("A21").NumberFormat = $ws.Range("A20").NumberFormat

# New row 22: date / question / URL (with hyperlink), mirroring the existing rows.
$ws.Range("A22").Value = $lisDate
$ws.Range("A22").NumberFormat = "mm-dd-yy"

$ws.Range("B22").Value = "LIS"

$lisUrl = "https://leetcode.com/problems/longest-increasing-subsequence/"
$ws.Hyperlinks.Add($ws.Range("C22"), $lisUrl, "", "", "LIS")

# Scroll the sheet view down so row 10 is at the top (matches authored view state).
$ws.Application.ActiveWindow.ScrollRow = 10
